$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 17

# Original layout:  A=Tiefe, B=v in cm/s(raw), C=v in m/s(formula =B/100, s="2"), D=Intens(raw; rows 12,16 use s="1")
# Target layout:    A=Tiefe, B=Tiefe in mm(NEW formula), C=v in cm/s(raw), D=v in m/s(formula =C/100, s="2"), E=Intens(raw; rows 12,16 use s="1")
# A new column is inserted at B; old B,C,D shift right to C,D,E respectively.
# Work right-to-left (E, D, C, B) so source data is not clobbered before being read,
# and explicitly clear/re-apply number formats so no stale styles leak across columns.

# --- Step 1: old D (Intens, raw numbers; rows 12 & 16 use "0.0") -> E ---
$ws.Range("E1").Formula = $ws.Range("D1").Formula
for ($r = 2; $r -le $lastRow; $r++) {
    $srcFmt = $ws.Cells.Item($r, 4).NumberFormat
    $ws.Cells.Item($r, 5).Formula = $ws.Cells.Item($r, 4).Formula
    $ws.Cells.Item($r, 5).ClearFormats()
    if ($srcFmt -ne "General") {
        $ws.Cells.Item($r, 5).NumberFormat = $srcFmt
    }
}

# --- Step 2: old C (v in m/s, formula =B/100, s="2" -> "0.000") -> D ---
$ws.Range("D1").Formula = $ws.Range("C1").Formula
$ws.Range("D2:D17").ClearFormats()
$ws.Range("D2:D17").Formula = "=C2/100"
$ws.Range("D2:D17").NumberFormat = "0.000"

# --- Step 3: old B (v in cm/s, raw numbers) -> C ---
$ws.Range("C1").Formula = $ws.Range("B1").Formula
for ($r = 2; $r -le $lastRow; $r++) {
    $srcFmt = $ws.Cells.Item($r, 2).NumberFormat
    $ws.Cells.Item($r, 3).Formula = $ws.Cells.Item($r, 2).Formula
    $ws.Cells.Item($r, 3).ClearFormats()
    if ($srcFmt -ne "General") {
        $ws.Cells.Item($r, 3).NumberFormat = $srcFmt
    }
}

# --- Step 4: new column B (Tiefe in mm) ---
$ws.Range("B1").Formula = "Tiefe in mm"
$ws.Range("B2:B17").ClearFormats()
$ws.Range("B2:B17").Formula = "=(A2-12.28)/4*6"

Write-Host "done"
